$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 2.891504666666667
$ws.Cells.Item(2,8).Value = 8.674514
$ws.Cells.Item(2,9).Value = 0.1213590456377548
$ws.Cells.Item(2,10).Value = 0.1213590456377548
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 22.906497
$ws.Cells.Item(2,14).Value = 68.719491
$ws.Cells.Item(2,15).Value = 0.9446038650914245
$ws.Cells.Item(2,16).Value = 0.9446038650914245
$ws.Cells.Item(2,17).Value = 66.234242972486
$ws.Cells.Item(2,18).Value = 596.1081867523741
$ws.Cells.Item(2,19).Value = 0.1146362235732297
$ws.Cells.Item(2,20).Value = 0.1146362235732297

# Row 3: ECs -> FAPs
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 2.891504666666667
$ws.Cells.Item(3,8).Value = 8.674514
$ws.Cells.Item(3,9).Value = 0.1213590456377548
$ws.Cells.Item(3,10).Value = 0.1213590456377548
$ws.Cells.Item(3,11).Value = 2
$ws.Cells.Item(3,12).Value = 0.6666666666666666
$ws.Cells.Item(3,13).Value = 0.1329193333333333
$ws.Cells.Item(3,14).Value = 0.3987579999999999
$ws.Cells.Item(3,15).Value = 0.005481244732096839
$ws.Cells.Item(3,16).Value = 0.005481244732096839
$ws.Cells.Item(3,17).Value = 0.3843368726235555
$ws.Cells.Item(3,18).Value = 3.459031853612
$ws.Cells.Item(3,19).Value = 0.0006651986295942433
$ws.Cells.Item(3,20).Value = 0.0006651986295942433

# Row 4: ECs -> sCs
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 2.891504666666667
$ws.Cells.Item(4,8).Value = 8.674514
$ws.Cells.Item(4,9).Value = 0.1213590456377548
$ws.Cells.Item(4,10).Value = 0.1213590456377548
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 1.210428333333333
$ws.Cells.Item(4,14).Value = 3.631285
$ws.Cells.Item(4,15).Value = 0.04991489017647865
$ws.Cells.Item(4,16).Value = 0.04991489017647865
$ws.Cells.Item(4,17).Value = 3.499959174498889
$ws.Cells.Item(4,18).Value = 31.49963257049
$ws.Cells.Item(4,19).Value = 0.00605762343493079
$ws.Cells.Item(4,20).Value = 0.006057623434930789

# Row 5: FAPs -> ECs
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 12.04042966666667
$ws.Cells.Item(5,8).Value = 36.121289
$ws.Cells.Item(5,9).Value = 0.505347637947847
$ws.Cells.Item(5,10).Value = 0.505347637947847
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 22.906497
$ws.Cells.Item(5,14).Value = 68.719491
$ws.Cells.Item(5,15).Value = 0.9446038650914245
$ws.Cells.Item(5,16).Value = 0.9446038650914245
$ws.Cells.Item(5,17).Value = 275.804066038211
$ws.Cells.Item(5,18).Value = 2482.236594343899
$ws.Cells.Item(5,19).Value = 0.4773533320203581
$ws.Cells.Item(5,20).Value = 0.4773533320203581

# Row 6: FAPs -> FAPs
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 12.04042966666667
$ws.Cells.Item(6,8).Value = 36.121289
$ws.Cells.Item(6,9).Value = 0.505347637947847
$ws.Cells.Item(6,10).Value = 0.505347637947847
$ws.Cells.Item(6,11).Value = 2
$ws.Cells.Item(6,12).Value = 0.6666666666666666
$ws.Cells.Item(6,13).Value = 0.1329193333333333
$ws.Cells.Item(6,14).Value = 0.3987579999999999
$ws.Cells.Item(6,15).Value = 0.005481244732096839
$ws.Cells.Item(6,16).Value = 0.005481244732096839
$ws.Cells.Item(6,17).Value = 1.600405884340222
$ws.Cells.Item(6,18).Value = 14.403652959062
$ws.Cells.Item(6,19).Value = 0.002769934078379217
$ws.Cells.Item(6,20).Value = 0.002769934078379217

# Row 7: FAPs -> sCs
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 12.04042966666667
$ws.Cells.Item(7,8).Value = 36.121289
$ws.Cells.Item(7,9).Value = 0.505347637947847
$ws.Cells.Item(7,10).Value = 0.505347637947847
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 1.210428333333333
$ws.Cells.Item(7,14).Value = 3.631285
$ws.Cells.Item(7,15).Value = 0.04991489017647865
$ws.Cells.Item(7,16).Value = 0.04991489017647865
$ws.Cells.Item(7,17).Value = 14.57407721404056
$ws.Cells.Item(7,18).Value = 131.166694926365
$ws.Cells.Item(7,19).Value = 0.02522437184910968
$ws.Cells.Item(7,20).Value = 0.02522437184910968

# Row 8: sCs -> ECs
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Mdk"
$ws.Cells.Item(8,3).Value = "Itga4"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 8.894099
$ws.Cells.Item(8,8).Value = 26.682297
$ws.Cells.Item(8,9).Value = 0.3732933164143983
$ws.Cells.Item(8,10).Value = 0.3732933164143982
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 22.906497
$ws.Cells.Item(8,14).Value = 68.719491
$ws.Cells.Item(8,15).Value = 0.9446038650914245
$ws.Cells.Item(8,16).Value = 0.9446038650914245
$ws.Cells.Item(8,17).Value = 203.732652061203
$ws.Cells.Item(8,18).Value = 1833.593868550827
$ws.Cells.Item(8,19).Value = 0.3526143094978367
$ws.Cells.Item(8,20).Value = 0.3526143094978366

# Row 9: sCs -> FAPs
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Mdk"
$ws.Cells.Item(9,3).Value = "Itga4"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 8.894099
$ws.Cells.Item(9,8).Value = 26.682297
$ws.Cells.Item(9,9).Value = 0.3732933164143983
$ws.Cells.Item(9,10).Value = 0.3732933164143982
$ws.Cells.Item(9,11).Value = 2
$ws.Cells.Item(9,12).Value = 0.6666666666666666
$ws.Cells.Item(9,13).Value = 0.1329193333333333
$ws.Cells.Item(9,14).Value = 0.3987579999999999
$ws.Cells.Item(9,15).Value = 0.005481244732096839
$ws.Cells.Item(9,16).Value = 0.005481244732096839
$ws.Cells.Item(9,17).Value = 1.182197709680667
$ws.Cells.Item(9,18).Value = 10.639779387126
$ws.Cells.Item(9,19).Value = 0.002046112024123379
$ws.Cells.Item(9,20).Value = 0.002046112024123379

# Row 10: sCs -> sCs
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Mdk"
$ws.Cells.Item(10,3).Value = "Itga4"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 8.894099
$ws.Cells.Item(10,8).Value = 26.682297
$ws.Cells.Item(10,9).Value = 0.3732933164143983
$ws.Cells.Item(10,10).Value = 0.3732933164143982
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 1.210428333333333
$ws.Cells.Item(10,14).Value = 3.631285
$ws.Cells.Item(10,15).Value = 0.04991489017647865
$ws.Cells.Item(10,16).Value = 0.04991489017647865
$ws.Cells.Item(10,17).Value = 10.76566942907167
$ws.Cells.Item(10,18).Value = 96.89102486164501
$ws.Cells.Item(10,19).Value = 0.01863289489243819
$ws.Cells.Item(10,20).Value = 0.01863289489243818
